$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Row 9: add a new journal entry (date + event), mirroring the style of row 8.
# Copy A8's formatting (date number format, alignment, etc.) down into A9 first,
# then set the actual date value (2020-04-09 = serial 43930).
$ws.Range("A8").Copy($ws.Range("A9"))
$ws.Range("A9").Value = 43930

# B9 gets the new event text, same style it already has.
$ws.Range("B9").Value = "Rendu de la documentation finale du projet"

# Move the active selection to B9, as in the saved workbook.
$ws.Range("B9").Select()
